$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add the "_Toc5728325" bookmark around the "Log" heading text.
#    (This also causes the pre-existing "_GoBack" bookmark to be renumbered
#    from id 0 to id 1, matching the target revision.)
# ---------------------------------------------------------------------------
$logParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Log`r") {
        $logParaIndex = $i
        break
    }
}
if ($logParaIndex -eq 0) {
    throw "Could not locate the 'Log' heading paragraph"
}
$logPara = $d.Paragraphs.Item($logParaIndex)
$logTextRange = $d.Range($logPara.Range.Start, $logPara.Range.Start + 3)
$d.Bookmarks.Add("_Toc5728325", $logTextRange)

# ---------------------------------------------------------------------------
# 2) The legacy VML text-box ("Text Box 138") and the INCLUDEPICTURE field
#    nested inside it are not reachable through the Shapes / Paragraphs
#    object model (they live in a story Word's automation surface here
#    does not enumerate), so those two edits are applied as raw OOXML
#    surgery through the Document.WordOpenXML round-trip property -
#    exactly like Word itself represents "Edit > Paste as RTF/XML" content.
# ---------------------------------------------------------------------------
$xml = $d.WordOpenXML

# 2a) Text Box 138 gets a new drawing id / z-index (Word re-mints these
#     whenever a shape bearing a field is touched).
$oldShapeAttrs = 'o:spid="_x0000_s1122" type="#_x0000_t202" style="position:absolute;margin-left:0;margin-top:0;width:134.85pt;height:302.4pt;z-index:251708416;'
$newShapeAttrs = 'o:spid="_x0000_s1134" type="#_x0000_t202" style="position:absolute;margin-left:0;margin-top:0;width:134.85pt;height:302.4pt;z-index:251714560;'
if ($xml.IndexOf($oldShapeAttrs) -lt 0) {
    throw "Shape attribute anchor not found"
}
$xml = $xml.Replace($oldShapeAttrs, $newShapeAttrs)

# 2b) Wrap the rendered INCLUDEPICTURE result (the w:pict run) in one more
#     nested field (begin / instrText / separate ... end), matching the
#     existing nesting pattern already present around it.
$rPr = "<w:rPr><w:b/><w:bCs/><w:sz w:val=`"26`"/><w:szCs w:val=`"26`"/></w:rPr>"
$instr = "<w:instrText xml:space=`"preserve`"> INCLUDEPICTURE  `"C:\\Users\\mikym\\Documents\\UniMi\\2S\\OGD\\Lama - OGDAIVG\\Lama\\lama.png`" \* MERGEFORMATINET </w:instrText>"
$pict = "<w:pict><v:shape id=`"_x0000_i1030`" type=`"#_x0000_t75`" style=`"width:56.25pt;height:57pt`"><v:imagedata r:id=`"rId8`" r:href=`"rId9`"/></v:shape></w:pict>"

$oldField = "<w:r>$rPr$pict</w:r><w:r>$rPr<w:fldChar w:fldCharType=`"end`"/></w:r>"
$newField = "<w:r>$rPr<w:fldChar w:fldCharType=`"begin`"/></w:r>" +
            "<w:r>$rPr$instr</w:r>" +
            "<w:r>$rPr<w:fldChar w:fldCharType=`"separate`"/></w:r>" +
            "<w:r>$rPr$pict</w:r>" +
            "<w:r>$rPr<w:fldChar w:fldCharType=`"end`"/></w:r>" +
            "<w:r>$rPr<w:fldChar w:fldCharType=`"end`"/></w:r>"

if ($xml.IndexOf($oldField) -lt 0) {
    throw "INCLUDEPICTURE field anchor not found"
}
$xml = $xml.Replace($oldField, $newField)

$d.WordOpenXML = $xml

Write-Output "Edit applied."
